# Update the Costing / Budget slide (slide 18) with the cost figures that
# were filled in for the "Development", "Infrastructure & Licensing" and
# "Operational" cost categories, and correct the "Initial Development &
# Setup Costs" range in the "Total Estimated Cost" summary textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)

# "Content Placeholder 3" (shape id 4) - the assumptions / cost breakdown body
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Paragraphs(5, 1).Runs(1, 1).Text = "Development Costs: `$31,500 - `$65,000"
$body.Paragraphs(9, 1).Runs(1, 1).Text = "Infrastructure & Licensing Costs: `$3,000 - `$12,500"
$body.Paragraphs(12, 1).Runs(1, 1).Text = "Operational Costs: `$26,000 - `$66,000"

# "Content Placeholder 3" (shape id 5) - the Total Estimated Cost textbox
$summary = $s.Shapes.Item(3).TextFrame.TextRange
$summary.Paragraphs(2, 1).Runs(1, 1).Text = "Initial Development & Setup Costs: `$35,000 - `$78,000"
